$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 09:52"

# 2. Rusia (row 11) - refreshed statistics, same rank position
$ws.Range("B11").Value = 114431
$ws.Range("C11").Value = 7933
$ws.Range("D11").Value = 13220
$ws.Range("E11").Value = 100042
$ws.Range("F11").Value = 2300
$ws.Range("G11").Value = 96
$ws.Range("H11").Value = 1169

# 3. Ucrania overtakes Corea del Sur (rows 38-39)
$ws.Range("A38").Value = "Ucrania"
$ws.Range("B38").Value = 10861
$ws.Range("C38").Value = 455
$ws.Range("D38").Value = 1413
$ws.Range("E38").Value = 9176
$ws.Range("F38").Value = 143
$ws.Range("G38").Value = 11
$ws.Range("H38").Value = 272

$ws.Range("A39").Value = "Corea del Sur"
$ws.Range("B39").Value = 10774
$ws.Range("C39").Value = 9
$ws.Range("D39").Value = 9072
$ws.Range("E39").Value = 1454
$ws.Range("F39").Value = 55
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 248

# 4. Lituania overtakes Eslovaquia (rows 85-86)
$ws.Range("A85").Value = "Lituania"
$ws.Range("B85").Value = 1399
$ws.Range("C85").Value = 14
$ws.Range("D85").Value = 594
$ws.Range("E85").Value = 760
$ws.Range("F85").Value = 17
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 45

$ws.Range("A86").Value = "Eslovaquia"
$ws.Range("B86").Value = 1396
$ws.Range("C86").Value = 0
$ws.Range("D86").Value = 524
$ws.Range("E86").Value = 849
$ws.Range("F86").Value = 8
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 23

# 5. Gambia jumps ahead of Santa Sede, Montserrat, Burundi and Seychelles (rows 201-205)
$ws.Range("A201").Value = "Gambia"
$ws.Range("B201").Value = 12
$ws.Range("C201").Value = 1
$ws.Range("D201").Value = 8
$ws.Range("E201").Value = 3
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 1

$ws.Range("A202").Value = "Santa Sede"
$ws.Range("B202").Value = 11
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 2
$ws.Range("E202").Value = 9
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

$ws.Range("A203").Value = "Montserrat"
$ws.Range("B203").Value = 11
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 2
$ws.Range("E203").Value = 8
$ws.Range("F203").Value = 1
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 1

$ws.Range("A204").Value = "Burundi"
$ws.Range("B204").Value = 11
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 4
$ws.Range("E204").Value = 6
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 1

$ws.Range("A205").Value = "Seychelles"
$ws.Range("B205").Value = 11
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 6
$ws.Range("E205").Value = 5
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0

# row 206 (Groenlandia) is unchanged

# 6. Comoras overtakes San Pedro y Miquelon (rows 217-218); figures are identical for both
$ws.Range("A217").Value = "Comoras"
$ws.Range("B217").Value = 1
$ws.Range("C217").Value = 0
$ws.Range("D217").Value = 0
$ws.Range("E217").Value = 1
$ws.Range("F217").Value = 0
$ws.Range("G217").Value = 0
$ws.Range("H217").Value = 0

$ws.Range("A218").Value = "San Pedro y Miquelon"
$ws.Range("B218").Value = 1
$ws.Range("C218").Value = 0
$ws.Range("D218").Value = 0
$ws.Range("E218").Value = 1
$ws.Range("F218").Value = 0
$ws.Range("G218").Value = 0
$ws.Range("H218").Value = 0
